$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 6
    3 = 3
    4 = 2
    5 = -2
    6 = -2
    7 = -2
    8 = 2
    9 = 2
    10 = -3
    12 = -1
    13 = 1
    14 = -7
    15 = 1
    16 = -3
    17 = -1
    18 = -1
    19 = -4
    20 = 2
    21 = -2
    22 = 5
    23 = 1
    24 = 4
    25 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
